$p = $ppt.ActivePresentation
$s = $p.Slides.Item(6)
$shape = $s.Shapes.Item(1)

$shape.Left = 49.999925
$shape.Width = 706.8964566929134

$shape.TextFrame.TextRange.Text = "Sort Results by Form Last Updated Date"
